$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# New row 19 values
$ws.Range("C19").Value = -101334.85
$ws.Range("E19").Formula = "=SUM(E18,I19,C19,D19)"
$ws.Range("F19").Value = "23-05-2014"

# Copy style from row 18 equivalents so formats match (E19 uses style of E18, F19 uses style of F18)
$ws.Range("E18").Copy()
$ws.Range("E19").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F18").Copy()
$ws.Range("F19").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("C20").Select()
